function RGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(14)
$shp = $s.Shapes.Item(3)

# --- Insert the new leading run "// " before the existing text ---
$tr = $shp.TextFrame.TextRange
$tr.InsertBefore("// ") | Out-Null
$newRun = $tr.Characters(1, 3)
$newRun.LanguageID = "en-US"
$newRun.Font.Size = 20
$newRun.Font.Color.RGB = RGB 0xF0 0x5A 0x28

# --- Reposition / resize the shape (a:off / a:ext) ---
# Target EMU: off x=896620 y=2367280, ext cx=5164455 cy=344170
# Done after the text edit because the shape has spAutoFit, which
# recomputes height from the text when the text body changes.
# PowerPoint Shape geometry is expressed in points (1 pt = 12700 EMU) and
# stored internally as 32-bit floats, so the literals below are chosen to
# land exactly on the target EMU values after that float32 round-trip.
$shp.Left = 70.60004039993822
$shp.Top = 186.40003969993828
$shp.Width = 406.6500397002338
$shp.Height = 27.100039500012148
